# Update "想去人数" (column F) figures on the "展览" and "全部类型" sheets
# to reflect the newer scrape snapshot (output generated at 456a3b4).

$wb = $excel.ActiveWorkbook

# Sheet "展览" (sheet1.xml): row -> new F value
$ws1 = $wb.Worksheets.Item("展览")
$ws1.Range("F3").Value  = 258
$ws1.Range("F5").Value  = 163
$ws1.Range("F6").Value  = 150
$ws1.Range("F8").Value  = 4603
$ws1.Range("F12").Value = 461
$ws1.Range("F14").Value = 15
$ws1.Range("F15").Value = 1326
$ws1.Range("F16").Value = 2626
$ws1.Range("F17").Value = 379
$ws1.Range("F19").Value = 55
$ws1.Range("F21").Value = 2294
$ws1.Range("F26").Value = 108

# Sheet "全部类型" (sheet4.xml): row -> new F value
$ws4 = $wb.Worksheets.Item("全部类型")
$ws4.Range("F3").Value  = 258
$ws4.Range("F5").Value  = 163
$ws4.Range("F6").Value  = 150
$ws4.Range("F9").Value  = 4603
$ws4.Range("F13").Value = 461
$ws4.Range("F15").Value = 15
$ws4.Range("F16").Value = 1326
$ws4.Range("F17").Value = 2626
$ws4.Range("F18").Value = 379
$ws4.Range("F20").Value = 55
$ws4.Range("F22").Value = 2294
$ws4.Range("F27").Value = 108
